# Update cryptos list data (price / volume change) per latest scrape run
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value is a plain decimal number need the source column's
# original Text formatting re-asserted first, otherwise COM auto-coerces
# the assignment to a Number and the "96.123.45"-style text values elsewhere
# in the same column would no longer match (column D is text, not numeric).
$textCells = @("D5", "D7", "D12", "D13", "D15", "D20", "D21", "D23", "D24", "D26", "D27", "D29", "D30", "D31", "D36", "D37", "D38", "D39", "D42", "D43", "D45", "D46", "D47", "D48", "D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "96.474.41"
$ws.Range("E2").Value = "  -0.60%  "
$ws.Range("D3").Value = "3.693.42"
$ws.Range("E3").Value = "  -0.40%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").Value = "236.18"
$ws.Range("E5").Value = "  -3.34%  "
$ws.Range("E6").Value = "  -1.31%  "
$ws.Range("D7").Value = "650.47"
$ws.Range("E7").Value = "  -1.56%  "
$ws.Range("E8").Value = "  -0.44%  "
$ws.Range("E9").Value = "  +0.01%  "
$ws.Range("E10").Value = "  -3.32%  "
$ws.Range("D11").Value = "3.690.96"
$ws.Range("E11").Value = "  -0.57%  "
$ws.Range("D12").Value = "0.0000311"
$ws.Range("E12").Value = "  +18.94%  "
$ws.Range("D13").Value = "44.23"
$ws.Range("E13").Value = "  -1.06%  "
$ws.Range("E14").Value = "  -0.22%  "
$ws.Range("D15").Value = "6.72"
$ws.Range("E15").Value = "  +2.51%  "
$ws.Range("D16").Value = "4.380.04"
$ws.Range("E16").Value = "  -0.33%  "
$ws.Range("D17").Value = "96.165.69"
$ws.Range("E17").Value = "  -1.00%  "
$ws.Range("E18").Value = "  +7.74%  "
$ws.Range("D19").Value = "3.687.62"
$ws.Range("E19").Value = "  +0.40%  "
$ws.Range("D20").Value = "13.12"
$ws.Range("E20").Value = "  +0.12%  "
$ws.Range("D21").Value = "18.64"
$ws.Range("E21").Value = "  +0.39%  "
$ws.Range("E22").Value = "  -5.95%  "
$ws.Range("D23").Value = "518.38"
$ws.Range("E23").Value = "  +0.62%  "
$ws.Range("D24").Value = "3.38"
$ws.Range("E24").Value = "  -2.52%  "
$ws.Range("E25").Value = "  +0.10%  "
$ws.Range("D26").Value = "6.92"
$ws.Range("E26").Value = "  -0.09%  "
$ws.Range("D27").Value = "100.55"
$ws.Range("E27").Value = "  -0.90%  "
$ws.Range("E28").Value = "  -0.57%  "
$ws.Range("D29").Value = "0.175"
$ws.Range("E29").Value = "  +3.43%  "
$ws.Range("D30").Value = "3.00"
$ws.Range("E30").Value = "  -1.23%  "
$ws.Range("D31").Value = "12.13"
$ws.Range("E31").Value = "  +0.51%  "
$ws.Range("E32").Value = "  -0.10%  "
$ws.Range("E33").Value = "  +6.61%  "
$ws.Range("E34").Value = "  -1.08%  "
$ws.Range("E35").Value = "  +0.41%  "
$ws.Range("B36").Value = "EthereumClassic"
$ws.Range("C36").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D36").Value = "32.08"
$ws.Range("E36").Value = "  -4.57%  "
$ws.Range("B37").Value = "Bittensor"
$ws.Range("C37").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D37").Value = "650.71"
$ws.Range("E37").Value = "  +5.06%  "
$ws.Range("D38").Value = "0.585"
$ws.Range("E38").Value = "  -0.98%  "
$ws.Range("D39").Value = "8.77"
$ws.Range("E39").Value = "  +0.17%  "
$ws.Range("E41").Value = "  +11.97%  "
$ws.Range("D42").Value = "2.06"
$ws.Range("E42").Value = "  +5.73%  "
$ws.Range("D43").Value = "40.42"
$ws.Range("E43").Value = "  -6.28%  "
$ws.Range("E44").Value = "  -0.39%  "
$ws.Range("D45").Value = "0.953"
$ws.Range("E45").Value = "  -2.29%  "
$ws.Range("D46").Value = "0.0450"
$ws.Range("E46").Value = "  +1.83%  "
$ws.Range("D47").Value = "0.428"
$ws.Range("E47").Value = "  +2.43%  "
$ws.Range("D48").Value = "23.59"
$ws.Range("E48").Value = "  -0.17%  "
$ws.Range("E49").Value = "  -1.73%  "
$ws.Range("E50").Value = "  -1.86%  "
$ws.Range("D51").Value = "3.53"
$ws.Range("E51").Value = "  +2.38%  "
